$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2684.0845
$ws.Range("I15").Value = 2684.0845
$ws.Range("K15").Value = 8052.2535
$ws.Range("M15").Value = -7883.2535
$ws.Range("H18").Value = 287.5
$ws.Range("I18").Value = 260
$ws.Range("J18").Value = 333.33334
$ws.Range("K18").Value = 260
$ws.Range("L18").Value = 333.33334
$ws.Range("M18").Value = 24
$ws.Range("N18").Value = -901.33334
$ws.Range("H32").Value = 1087.1333
$ws.Range("I32").Value = 443.33334
$ws.Range("J32").Value = 1158.6666
$ws.Range("K32").Value = 443.33334
$ws.Range("L32").Value = 1158.6666
$ws.Range("M32").Value = -117.33334
$ws.Range("N32").Value = -1810.6666
$ws.Range("H54").Value = 6500
$ws.Range("I54").Value = 3000
$ws.Range("J54").Value = 10000
$ws.Range("K54").Value = 3000
$ws.Range("L54").Value = 10000
$ws.Range("M54").Value = -2514
$ws.Range("N54").Value = -10972
$ws.Range("H137").Value = 8930271
$ws.Range("I137").Value = 15153392
$ws.Range("J137").Value = 1445.7391
$ws.Range("K137").Value = 45460176
$ws.Range("L137").Value = 4337.2173
$ws.Range("M137").Value = -45457626
$ws.Range("N137").Value = -9437.2173
$ws.Range("H140").Value = 76195
$ws.Range("J140").Value = 76195
$ws.Range("L140").Value = 76195
$ws.Range("N140").Value = -86555
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 3697.7144
$ws.Range("I25").Value = 2076.8
$ws.Range("J25").Value = 7750
$ws.Range("K25").Value = 2076.8
$ws.Range("L25").Value = 7750
$ws.Range("M25").Value = -1674.8
$ws.Range("N25").Value = -8554
$ws.Range("H45").Value = 3852
$ws.Range("I45").Value = 2376.5
$ws.Range("J45").Value = 4737.3
$ws.Range("K45").Value = 2376.5
$ws.Range("L45").Value = 4737.3
$ws.Range("M45").Value = -1999.5
$ws.Range("N45").Value = -5491.3
$ws.Range("H57").Value = 27000
$ws.Range("I57").Value = 27000
$ws.Range("K57").Value = 27000
$ws.Range("M57").Value = -26516
$ws.Range("H61").Value = 2111.0698
$ws.Range("I61").Value = 1357.6842
$ws.Range("J61").Value = 2707.5
$ws.Range("K61").Value = 1357.6842
$ws.Range("L61").Value = 2707.5
$ws.Range("M61").Value = -1145.6842
$ws.Range("N61").Value = -3131.5
$ws.Range("H133").Value = 35151.832
$ws.Range("J133").Value = 35151.832
$ws.Range("L133").Value = 35151.832
$ws.Range("N133").Value = -40211.832
$ws.Range("H136").Value = 2111.0698
$ws.Range("I136").Value = 1357.6842
$ws.Range("J136").Value = 2707.5
$ws.Range("K136").Value = 4073.0526
$ws.Range("L136").Value = 8122.5
$ws.Range("M136").Value = -1523.0526
$ws.Range("N136").Value = -13222.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 765.45
$ws.Range("I22").Value = 881.8125
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 881.8125
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -708.8125
$ws.Range("N22").Value = -646
$ws.Range("H37").Value = 1184.3334
$ws.Range("I37").Value = 415.4
$ws.Range("J37").Value = 5029
$ws.Range("K37").Value = 415.4
$ws.Range("L37").Value = 5029
$ws.Range("M37").Value = -278.4
$ws.Range("N37").Value = -5303
$ws.Range("H105").Value = 3799.2307
$ws.Range("I105").Value = 1963.75
$ws.Range("J105").Value = 6736
$ws.Range("K105").Value = 1963.75
$ws.Range("L105").Value = 6736
$ws.Range("M105").Value = -216.75
$ws.Range("N105").Value = -10230
$ws.Range("H128").Value = 2618
$ws.Range("I128").Value = 2618
$ws.Range("K128").Value = 7854
$ws.Range("M128").Value = -5364
$ws.Range("H134").Value = 5204.523
$ws.Range("I134").Value = 2760.4707
$ws.Range("J134").Value = 6743.3706
$ws.Range("K134").Value = 8281.4121
$ws.Range("L134").Value = 20230.1118
$ws.Range("M134").Value = -5746.4121
$ws.Range("N134").Value = -25300.1118
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H33").Value = 3000
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2621
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H132").Value = 3559.3333
$ws.Range("I132").Value = 2154.6667
$ws.Range("J132").Value = 5666.3335
$ws.Range("K132").Value = 6464.000100000001
$ws.Range("L132").Value = 16999.0005
$ws.Range("M132").Value = -3934.000100000001
$ws.Range("N132").Value = -22059.0005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 165.6
$ws.Range("I4").Value = 128.44444
$ws.Range("K4").Value = 385.33332
$ws.Range("M4").Value = -273.33332
$ws.Range("H63").Value = 2547.111
$ws.Range("I63").Value = 904.8
$ws.Range("J63").Value = 4600
$ws.Range("K63").Value = 2714.4
$ws.Range("L63").Value = 13800
$ws.Range("M63").Value = -1965.4
$ws.Range("N63").Value = -15298
$ws.Range("H66").Value = 2547.111
$ws.Range("I66").Value = 904.8
$ws.Range("J66").Value = 4600
$ws.Range("K66").Value = 8143.2
$ws.Range("L66").Value = 41400
$ws.Range("M66").Value = -4399.2
$ws.Range("N66").Value = -48888
$ws.Range("H68").Value = 555239.75
$ws.Range("I68").Value = 794.4400000000001
$ws.Range("J68").Value = 798417.5
$ws.Range("K68").Value = 2383.32
$ws.Range("L68").Value = 2395252.5
$ws.Range("M68").Value = -1572.32
$ws.Range("N68").Value = -2396874.5
$ws.Range("H71").Value = 555239.75
$ws.Range("I71").Value = 794.4400000000001
$ws.Range("J71").Value = 798417.5
$ws.Range("K71").Value = 7149.960000000001
$ws.Range("L71").Value = 7185757.5
$ws.Range("M71").Value = -3093.960000000001
$ws.Range("N71").Value = -7193869.5
$ws.Range("H107").Value = 23207220
$ws.Range("J107").Value = 1504592
$ws.Range("L107").Value = 4513776
$ws.Range("N107").Value = -4517616
$ws.Range("H114").Value = 1214.5
$ws.Range("I114").Value = 328
$ws.Range("J114").Value = 1657.75
$ws.Range("K114").Value = 984
$ws.Range("L114").Value = 4973.25
$ws.Range("M114").Value = 2270
$ws.Range("N114").Value = -11481.25
$ws.Range("H117").Value = 17857564
$ws.Range("I117").Value = 317
$ws.Range("J117").Value = 35714812
$ws.Range("K117").Value = 951
$ws.Range("L117").Value = 107144436
$ws.Range("M117").Value = 2491
$ws.Range("N117").Value = -107151320
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 5000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4888
$ws.Range("N5").ClearContents()
$ws.Range("H132").Value = 2832.2307
$ws.Range("I132").Value = 2438.2727
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 7314.8181
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4784.8181
$ws.Range("N132").Value = -20057
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H46").Value = 834357.9399999999
$ws.Range("I46").Value = 843.7143
$ws.Range("K46").Value = 843.7143
$ws.Range("M46").Value = -655.7143
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 260035000
$ws.Range("I2").Value = 40000000
$ws.Range("J2").Value = 333380000
$ws.Range("K2").Value = 40000000
$ws.Range("L2").Value = 333380000
$ws.Range("M2").Value = -39999888
$ws.Range("N2").Value = -333380224
$ws.Range("H26").Value = 70014
$ws.Range("J26").Value = 70014
$ws.Range("L26").Value = 70014
$ws.Range("N26").Value = -70600
$ws.Range("H28").Value = 26006.334
$ws.Range("J28").Value = 26006.334
$ws.Range("L28").Value = 26006.334
$ws.Range("N28").Value = -26702.334
$ws.Range("H47").Value = 6000
$ws.Range("J47").Value = 6000
$ws.Range("L47").Value = 6000
$ws.Range("N47").Value = -7144
$ws.Range("H54").Value = 7666.6665
$ws.Range("J54").Value = 7666.6665
$ws.Range("L54").Value = 7666.6665
$ws.Range("N54").Value = -8706.666499999999
$ws.Range("H122").Value = 7369.593
$ws.Range("I122").Value = 8606
$ws.Range("J122").Value = 4433.125
$ws.Range("K122").Value = 25818
$ws.Range("L122").Value = 13299.375
$ws.Range("M122").Value = -23368
$ws.Range("N122").Value = -18199.375